$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking price cells as Text so they keep their exact
# string representation (matching source data) instead of being parsed as numbers.
foreach ($addr in @("D5","D6","D8","D10","D11","D12","D14","D17","D19","D20","D21","D23","D24","D25","D27","D28","D29","D31","D32","D34","D35","D36","D37","D38","D40","D42","D47","D49","D50","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin data values
$ws.Range("D2").Value = "47.887.21"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "2.481.25"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "317.04"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").Value = "104.72"
$ws.Range("E6").Value = "  -3.88%  "
$ws.Range("E7").Value = "  -2.58%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -3.32%  "
$ws.Range("D10").Value = "38.82"
$ws.Range("E10").Value = "  -3.91%  "
$ws.Range("D11").Value = "20.35"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "0.0798"
$ws.Range("E12").Value = "  -2.99%  "
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "7.02"
$ws.Range("E14").Value = "  -3.26%  "
$ws.Range("D15").Value = "2.870.90"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "2.491.39"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").Value = "0.822"
$ws.Range("E17").Value = "  -3.72%  "
$ws.Range("D18").Value = "47.838.93"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "2.95"
$ws.Range("E19").Value = "  +8.82%  "
$ws.Range("D20").Value = "12.65"
$ws.Range("E20").Value = "  -4.32%  "
$ws.Range("D21").Value = "6.51"
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("D22").Value = "0.0₃0925"
$ws.Range("D23").Value = "278.94"
$ws.Range("E23").Value = "  +5.57%  "
$ws.Range("D24").Value = "70.76"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "25.61"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").Value = "9.56"
$ws.Range("E29").Value = "  -5.38%  "
$ws.Range("E30").Value = "  -3.99%  "
$ws.Range("D31").Value = "34.55"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").Value = "49.20"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "18.86"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "5.24"
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("D36").Value = "0.0767"
$ws.Range("E36").Value = "  -2.76%  "
$ws.Range("D37").Value = "1.93"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("D38").Value = "4.51"
$ws.Range("E38").Value = "  -4.14%  "
$ws.Range("E39").Value = "  -4.47%  "
$ws.Range("D40").Value = "122.16"
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").Value = "22.03"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "1.990.29"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").Value = "1.89"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("D49").Value = "8.91"
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("D50").Value = "5.12"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").Value = "78.93"
$ws.Range("E51").Value = "  +0.16%  "
